$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.680.94"
$ws.Range("E2").Value = "  +1.59%  "
$ws.Range("D3").Value = "2.298.37"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'323.29"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "'104.30"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").Value = "'40.15"
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "'8.42"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "'0.975"
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("D15").Value = "'15.28"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "2.645.22"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "2.281.98"
$ws.Range("E17").Value = "  -1.61%  "
$ws.Range("D18").Value = "42.593.30"
$ws.Range("E18").Value = "  +1.11%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  -3.18%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'13.65"
$ws.Range("E21").Value = "  +36.67%  "
$ws.Range("D22").Value = "'73.37"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'3.60"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'269.95"
$ws.Range("E24").Value = "  -7.11%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("E26").Value = "  +0.46%  "
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'38.61"
$ws.Range("E29").Value = "  +9.27%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.60"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").Value = "'165.68"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "'6.19"
$ws.Range("E32").Value = "  +4.86%  "
$ws.Range("D33").Value = "'0.0884"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("E36").Value = "  -13.98%  "
$ws.Range("D37").Value = "'4.64"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "'3.72"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").Value = "'2.75"
$ws.Range("E40").Value = "  -3.95%  "
$ws.Range("D41").Value = "'1.54"
$ws.Range("E41").Value = "  +4.17%  "
$ws.Range("D42").Value = "'69.63"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.226"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").Value = "'93.13"
$ws.Range("E45").Value = "  -9.82%  "
$ws.Range("D46").Value = "'12.41"
$ws.Range("E46").Value = "  +2.63%  "
$ws.Range("D47").Value = "'81.95"
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("D48").Value = "'113.67"
$ws.Range("E48").Value = "  -3.50%  "
$ws.Range("D49").Value = "'8.95"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").Value = "1.602.46"
$ws.Range("E51").Value = "  +2.73%  "

# Reset number format on cells that were forced to text so no stray
# text-format style is left on them (keeps styles.xml unchanged).
$resetRefs = @("D5","D6","D9","D10","D11","D12","D14","D15","D19","D21","D22","D23","D24","D25","D27","D29","D30","D31","D32","D33","D35","D37","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49")
foreach ($r in $resetRefs) {
    $ws.Range($r).Style = "Normal"
}
